$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.411.46"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "'3.686.06"
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'682.54"
$ws.Range("E5").Value = "  -2.63%  "
$ws.Range("D6").Value = "'162.52"
$ws.Range("E6").Value = "  -3.96%  "
$ws.Range("D7").Value = "'3.685.05"
$ws.Range("E7").Value = "  -2.88%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.500"
$ws.Range("E9").Value = "  -4.00%  "
$ws.Range("E10").Value = "  -7.39%  "
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("E13").Value = "  -4.12%  "
$ws.Range("D14").Value = "'33.50"
$ws.Range("E14").Value = "  -5.93%  "
$ws.Range("D15").Value = "'4.308.95"
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("D16").Value = "'3.689.34"
$ws.Range("E16").Value = "  -2.43%  "
$ws.Range("D17").Value = "'69.430.26"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  -6.01%  "
$ws.Range("E20").Value = "  -6.46%  "
$ws.Range("D21").Value = "'481.35"
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("D22").Value = "'9.89"
$ws.Range("E22").Value = "  -7.18%  "
$ws.Range("D23").Value = "'0.667"
$ws.Range("E23").Value = "  -7.51%  "
$ws.Range("D24").Value = "'80.30"
$ws.Range("E24").Value = "  -4.52%  "
$ws.Range("D25").Value = "'3.832.22"
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("E26").Value = "  -8.28%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'11.51"
$ws.Range("E28").Value = "  -3.94%  "
$ws.Range("D29").Value = "'9.59"
$ws.Range("E29").Value = "  -6.31%  "
$ws.Range("E30").Value = "  -8.32%  "
$ws.Range("E31").Value = "  -8.94%  "
$ws.Range("D32").Value = "'2.09"
$ws.Range("E32").Value = "  -7.86%  "
$ws.Range("D33").Value = "'6.85"
$ws.Range("E33").Value = "  -6.03%  "
$ws.Range("D34").Value = "'27.11"
$ws.Range("E34").Value = "  -6.22%  "
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "'3.653.93"
$ws.Range("E37").Value = "  -3.05%  "
$ws.Range("D38").Value = "'8.51"
$ws.Range("E38").Value = "  -5.60%  "
$ws.Range("E39").Value = "  +7.34%  "
$ws.Range("D40").Value = "'0.0936"
$ws.Range("E40").Value = "  -7.16%  "
$ws.Range("E41").Value = "  -3.77%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E44").Value = "  -6.22%  "
$ws.Range("D45").Value = "'160.20"
$ws.Range("E45").Value = "  -3.85%  "
$ws.Range("D46").Value = "'48.35"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "'2.85"
$ws.Range("E47").Value = "  -11.04%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'30.11"
$ws.Range("E48").Value = "  +7.21%  "
$ws.Range("D49").Value = "'0.000290"
$ws.Range("E49").Value = "  -7.72%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'1.36"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("D51").Value = "'395.24"
$ws.Range("E51").Value = "  -6.10%  "
